$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new vendor keyword entry in the next available row
$ws.Range("A12").Value = "wpnb_pto_new_users_add()"

# Column A no longer auto-fits; widen it to comfortably fit the longer keyword
$ws.Columns.Item(1).ColumnWidth = 24.53

# Move / update the selection to match where the user ended up after editing
$ws.Range("F17:F18").Select()
